$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 114
$ws.Range("B2").Value = "Anantara"
$ws.Range("C2").Value = "Test1"
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "360"
